$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" row (row 2) query text loses its trailing Cohort column
# (`coalesce(co.cohort_description, '') AS `Cohort``) from the RETURN clause.
$newCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nMATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nWHERE diag.stage_of_disease IN ['IVb']`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newCasesQuery

# Restore the default top-left/selection view state (scrolled view reset, selection moves to B2)
$ws.Range("B2").Select()

$wb.Save()
